# Generate Report for Handoff
#
# A new handoff cycle ran for the e2e markdown file: its source filename's
# GUID changed (a new source doc) and the generated xliff content hash
# changed too. The handback side has not happened yet for this cycle, so
# the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on both locale sheets get cleared back to "not yet
# available" (empty string / the zero DateTime).

$wb = $excel.ActiveWorkbook

$oldGuid = "17ae3688-f601-4fda-9ec4-75a4b1907617"
$newGuid = "d348dcc0-43b2-40cd-a538-ea3d30f62eb2"

$oldHash = "e388192d03ef0a6dd716f17b66d359b14047db4d"
$newHash = "84cec70876c16471171928c528108ba653fbc10a"

$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = ($newGuid + ".md")
$wsOverview.Range("B2").Value = ("e2e\" + $newGuid + ".md")
foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = ("e2e\" + $newGuid + ".md")
    }
}
$wsOverview.Range("G2").Value = "2016-09-01 23:05:10"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = ($newGuid + ".md")
foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = ($newGuid + ".md")
    }
}

$wsZh.Range("G2").Value = ($newGuid + "." + $newHash + ".zh-cn.xlf")
$wsZh.Range("H2").Value = "2016-09-01 23:04:58"

# This locale's target/handback has not happened for the new cycle yet;
# drop the stale target/handback file links + restore the hyperlink-free,
# default-styled look of an empty cell.
foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
# "'" forces a literal (Text) empty string instead of clearing the cell to
# blank, matching the pre-existing empty-string cells elsewhere in the row.
$wsZh.Range("I2").Value = "'"
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = "'"
$wsZh.Range("J2").Style = "Normal"
$wsZh.Range("K2").Value = $zeroDate

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = ($newGuid + ".md")
foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = ($newGuid + ".md")
    }
}

$wsDe.Range("G2").Value = ($newGuid + "." + $newHash + ".de-de.xlf")
$wsDe.Range("H2").Value = "2016-09-01 23:05:10"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
$wsDe.Range("I2").Value = "'"
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = "'"
$wsDe.Range("J2").Style = "Normal"
$wsDe.Range("K2").Value = $zeroDate
